$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width: 19.7109375 -> 23 ---
$ws.Columns("B").ColumnWidth = 22.14

# --- Write new word-pair values for rows 57-75 ---
# Written in the exact order that reproduces the original author's
# shared-string insertion sequence (indices 111-147).
$ws.Range("B57").Value = "chasing"
$ws.Range("C57").Value = "पीछा"
$ws.Range("B58").Value = "earnest"
$ws.Range("C58").Value = "ईमानदार"
$ws.Range("B59").Value = "eerie"
$ws.Range("C59").Value = "भयानक"
$ws.Range("C60").Value = "टालना"
$ws.Range("B60").Value = "elude"
$ws.Range("C61").Value = "शर्मनाक"
$ws.Range("B61").Value = "embarrassing"
$ws.Range("C62").Value = "काट - छाँट करना"
$ws.Range("B62").Value = "emend"
$ws.Range("C63").Value = "ज़ोर देना"
$ws.Range("B63").Value = "emphasise"
$ws.Range("B64").Value = "endeavour"
$ws.Range("C64").Value = "प्रयास"
$ws.Range("C65").Value = "महापाप"
$ws.Range("B65").Value = "enormity"
$ws.Range("B66").Value = "enquiry(british)"
$ws.Range("B67").Value = "Inquiry(American)"
$ws.Range("C66").Value = "both are correct"
$ws.Range("C68").Value = "ठीक कर लेना"
$ws.Range("B68").Value = "insure(arrange for financial)"
$ws.Range("B69").Value = "ensure"
$ws.Range("C69").Value = "सुनिश्चित करना"
$ws.Range("C70").Value = "अंत में"
$ws.Range("B70").Value = "eventually"
$ws.Range("B71").Value = "EXHAUSTED"
$ws.Range("C71").Value = "थका"
$ws.Range("C72").Value = "चिल्लाना"
$ws.Range("B72").Value = "exclaim"
$ws.Range("B73").Value = "explicit"
$ws.Range("C73").Value = "स्पष्ट"
$ws.Range("C74").Value = "अस्पष्ट"
$ws.Range("B74").Value = "implicit"
$ws.Range("C75").Value = "फिजूलखर्ची"
$ws.Range("B75").Value = "extravagance"

# --- Re-apply the existing "C" style (index 2: Arial 10pt black) to the ---
# --- newly-created cells C61:C65 (rows 57-60 already carried s="2"). ---
$ws.Range("C1").Copy()
$ws.Range("C61:C65").PasteSpecial(-4122)

# --- Row 67 column C: blank cell that only carries an alignment flag. ---
# This is cell-style index 3 in the target: default font + applyAlignment.
$ws.Range("C67").WrapText = $false

# --- Rows 66,68-75 column C: Arial(10pt,black) font + applyAlignment flag. ---
# This becomes cell-style index 4 in the target.
$ws.Range("C1").Copy()
$ws.Range("C66").PasteSpecial(-4122)
$ws.Range("C66").WrapText = $false

$ws.Range("C1").Copy()
$ws.Range("C68:C75").PasteSpecial(-4122)
$ws.Range("C68:C75").WrapText = $false

# --- Update selection to match the final active cell ---
$ws.Range("B75").Select()
